$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.071.59"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.645.25"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5047"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.012"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06441"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "1.648.43"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "1.873.20"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5461"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "0.0₅7935"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "26.085.18"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "204.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.313"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.971"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.013"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.938"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.752"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05066"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.244"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.264"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.199"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.545"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.349"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8980"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.621"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5644"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").Value = "1.150.36"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01575"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.577"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.012"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.677"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "1.784.05"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("D47").Value = "0.0₈112"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4543"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.011"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("E51").Value = "  -0.94%  "
